$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1146
$ws.Range("E2").Value = 99
$ws.Range("F2").Value = 99
$ws.Range("G2").Value = 104
$ws.Range("H2").Value = 76
$ws.Range("I2").Value = 76
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1615
$ws.Range("L2").Value = 414
$ws.Range("M2").Value = 1200
$ws.Range("N2").Value = 1195
$ws.Range("O2").Value = 5
$ws.Range("P2").Value = 68
$ws.Range("Q2").Value = 91
$ws.Range("R2").Value = 52
$ws.Range("S2").Value = -24
$ws.Range("T2").Value = 19
$ws.Range("U2").Value = 72
$ws.Range("W2").Value = 8.65
$ws.Range("X2").Value = 6.63
$ws.Range("Y2").Value = 6.48
$ws.Range("Z2").Value = 4.83
$ws.Range("AA2").Value = 34.53
$ws.Range("AB2").Value = 1786.28
$ws.Range("AC2").Value = 556
$ws.Range("AD2").Value = 11.6
$ws.Range("AE2").Value = 10304
$ws.Range("AF2").Value = 0.63
$ws.Range("AG2").Value = 200
$ws.Range("AH2").Value = 3.1
$ws.Range("AI2").Value = 30.64
$ws.Range("AJ2").Value = 13617577

# Row 3
$ws.Range("D3").Value = 1158
$ws.Range("E3").Value = 84
$ws.Range("F3").Value = 84
$ws.Range("G3").Value = 72
$ws.Range("H3").Value = 57
$ws.Range("I3").Value = 57
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1701
$ws.Range("L3").Value = 467
$ws.Range("M3").Value = 1234
$ws.Range("N3").Value = 1229
$ws.Range("O3").Value = 5
$ws.Range("P3").Value = 68
$ws.Range("Q3").Value = 29
$ws.Range("R3").Value = -25
$ws.Range("S3").Value = -24
$ws.Range("T3").Value = 29
$ws.Range("U3").Value = 0
$ws.Range("W3").Value = 7.29
$ws.Range("X3").Value = 4.94
$ws.Range("Y3").Value = 4.7
$ws.Range("Z3").Value = 3.45
$ws.Range("AA3").Value = 37.84
$ws.Range("AB3").Value = 1835.93
$ws.Range("AC3").Value = 419
$ws.Range("AD3").Value = 15.57
$ws.Range("AE3").Value = 10595
$ws.Range("AF3").Value = 0.62
$ws.Range("AG3").Value = 220
$ws.Range("AH3").Value = 3.37
$ws.Range("AI3").Value = 44.77
$ws.Range("AJ3").Value = 13617577

# Row 4
$ws.Range("D4").Value = 1158
$ws.Range("E4").Value = 61
$ws.Range("F4").Value = 61
$ws.Range("G4").Value = 67
$ws.Range("H4").Value = 59
$ws.Range("I4").Value = 59
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1698
$ws.Range("L4").Value = 431
$ws.Range("M4").Value = 1267
$ws.Range("N4").Value = 1262
$ws.Range("O4").Value = 5
$ws.Range("P4").Value = 68
$ws.Range("Q4").Value = 22
$ws.Range("R4").Value = -57
$ws.Range("S4").Value = -26
$ws.Range("T4").Value = 72
$ws.Range("U4").Value = -50
$ws.Range("W4").Value = 5.25
$ws.Range("X4").Value = 5.07
$ws.Range("Y4").Value = 4.71
$ws.Range("Z4").Value = 3.46
$ws.Range("AA4").Value = 33.99
$ws.Range("AB4").Value = 1884.67
$ws.Range("AC4").Value = 431
$ws.Range("AD4").Value = 16.96
$ws.Range("AE4").Value = 10881
$ws.Range("AF4").Value = 0.67
$ws.Range("AG4").Value = 250
$ws.Range("AH4").Value = 3.42
$ws.Range("AI4").Value = 49.41
$ws.Range("AJ4").Value = 13617577

# Row 5
$ws.Range("D5").Value = 1202
$ws.Range("E5").Value = 67
$ws.Range("F5").Value = 67
$ws.Range("G5").Value = 82
$ws.Range("H5").Value = 59
$ws.Range("I5").Value = 59
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1799
$ws.Range("L5").Value = 501
$ws.Range("M5").Value = 1297
$ws.Range("N5").Value = 1293
$ws.Range("O5").Value = 5
$ws.Range("P5").Value = 68
$ws.Range("Q5").Value = 93
$ws.Range("R5").Value = 108
$ws.Range("S5").Value = -29
$ws.Range("T5").Value = 22
$ws.Range("U5").Value = 71
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 5.54
$ws.Range("X5").Value = 4.92
$ws.Range("Y5").Value = 4.63
$ws.Range("Z5").Value = 3.38
$ws.Range("AA5").Value = 38.65
$ws.Range("AB5").Value = 1929.02
$ws.Range("AC5").Value = 435
$ws.Range("AD5").Value = 16.1
$ws.Range("AE5").Value = 11141
$ws.Range("AF5").Value = 0.63
$ws.Range("AG5").Value = 250
$ws.Range("AH5").Value = 3.57
$ws.Range("AI5").Value = 49
$ws.Range("AJ5").Value = 13617577

# Row 6
$ws.Range("D6").Value = 1211
$ws.Range("E6").Value = 68
$ws.Range("F6").Value = 68
$ws.Range("G6").Value = 75
$ws.Range("H6").Value = 62
$ws.Range("I6").Value = 61
$ws.Range("K6").Value = 1857
$ws.Range("L6").Value = 527
$ws.Range("M6").Value = 1330
$ws.Range("N6").Value = 1325
$ws.Range("P6").Value = 68
$ws.Range("Q6").Value = 27
$ws.Range("R6").Value = -30
$ws.Range("S6").Value = -30
$ws.Range("T6").Value = 34
$ws.Range("U6").Value = -8
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 5.64
$ws.Range("X6").Value = 5.08
$ws.Range("Y6").Value = 4.69
$ws.Range("Z6").Value = 3.37
$ws.Range("AA6").Value = 39.65
$ws.Range("AB6").Value = 1976.48
$ws.Range("AC6").Value = 450
$ws.Range("AD6").Value = 14.66
$ws.Range("AE6").Value = 11420
$ws.Range("AF6").Value = 0.58
$ws.Range("AG6").Value = 250
$ws.Range("AH6").Value = 3.79
$ws.Range("AI6").Value = 47.3
$ws.Range("AJ6").Value = 13617577

# Clear V column for rows 2-4 (column removed entirely)
$ws.Range("V2:V4").ClearContents()

# Clear rows 7-9 (only A,B,C retained)
$ws.Range("D7:AJ9").ClearContents()
